$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.150.53"
$ws.Range("E2").Value = "  +0.38%  "
$ws.Range("D3").Value = "1.675.14"
$ws.Range("E3").Value = "  -0.26%  "
$ws.Range("E4").Value = "  +0.09%  "
$ws.Range("D5").Value = "'214.35"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.68%  "
$ws.Range("E6").Value = "  +0.05%  "
$ws.Range("E7").Value = "  +0.06%  "
$ws.Range("D8").Value = "'22.77"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +6.33%  "
$ws.Range("E9").Value = "  +3.16%  "
$ws.Range("E10").Value = "  -0.32%  "
$ws.Range("E11").Value = "  +0.01%  "
$ws.Range("D12").Value = "1.912.83"
$ws.Range("E12").Value = "  -0.21%  "
$ws.Range("D13").Value = "1.681.88"
$ws.Range("E13").Value = "  +0.85%  "
$ws.Range("E14").Value = "  +2.51%  "
$ws.Range("D15").Value = "'0.559"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +4.95%  "
$ws.Range("D16").Value = "'66.56"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.35%  "
$ws.Range("D17").Value = "27.122.03"
$ws.Range("E17").Value = "  +0.28%  "
$ws.Range("D18").Value = "'235.04"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.62%  "
$ws.Range("D19").Value = "'7.86"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -3.66%  "
$ws.Range("E20").Value = "  +0.54%  "
$ws.Range("E21").Value = "  +0.12%  "
$ws.Range("D22").Value = "'4.54"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +1.79%  "
$ws.Range("D23").Value = "'9.55"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +3.15%  "
$ws.Range("D24").Value = "'2.09"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -1.89%  "
$ws.Range("D25").Value = "'148.05"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.53%  "
$ws.Range("D26").Value = "'7.46"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +2.54%  "
$ws.Range("E27").Value = "  -0.62%  "
$ws.Range("E28").Value = "  -0.31%  "
$ws.Range("E29").Value = "  +0.19%  "
$ws.Range("D30").Value = "'0.0501"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.75%  "
$ws.Range("E31").Value = "  -0.53%  "
$ws.Range("E32").Value = "  +0.18%  "
$ws.Range("D33").Value = "1.540.51"
$ws.Range("E33").Value = "  -0.25%  "
$ws.Range("E34").Value = "  +1.31%  "
$ws.Range("E35").Value = "  -3.83%  "
$ws.Range("D36").Value = "'0.607"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +3.11%  "
$ws.Range("D37").Value = "'0.944"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +3.30%  "
$ws.Range("E38").Value = "  +0.02%  "
$ws.Range("E39").Value = "  -0.87%  "
$ws.Range("E40").Value = "  +2.32%  "
$ws.Range("D41").Value = "'69.84"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +2.93%  "
$ws.Range("E42").Value = "  +4.56%  "
$ws.Range("E43").Value = "  +0.06%  "
$ws.Range("E44").Value = "  -0.36%  "
$ws.Range("D45").Value = "1.822.10"
$ws.Range("E45").Value = "  +0.01%  "
$ws.Range("D46").Value = "'0.781"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.22%  "
$ws.Range("D47").Value = "'89.71"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.86%  "
$ws.Range("E48").Value = "  +6.67%  "
$ws.Range("E49").Value = "  +3.76%  "
$ws.Range("D50").Value = "'8.26"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +2.56%  "
$ws.Range("E51").Value = "  -0.19%  "
